# Updates the cryptos price list (Price column D, Volume(1h) change column E)
# with freshly scraped values, per the GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell values scraped for this run: cell reference -> new text.
$updates = @{
    'D2' = '62.248.69'
    'E2' = '  +1.35%  '
    'D3' = '2.420.20'
    'E3' = '  +1.74%  '
    'E4' = '  -0.03%  '
    'D5' = '561.68'
    'E5' = '  +1.59%  '
    'D6' = '143.70'
    'E6' = '  +2.86%  '
    'E7' = '  +0.05%  '
    'E8' = '  +1.43%  '
    'D9' = '2.417.29'
    'E9' = '  +1.54%  '
    'E10' = '  +1.01%  '
    'E11' = '  -2.10%  '
    'D12' = '5.37'
    'E12' = '  +0.14%  '
    'E13' = '  +0.49%  '
    'D14' = '25.87'
    'E14' = '  +1.24%  '
    'D15' = '0.0000176'
    'E15' = '  +2.07%  '
    'D16' = '2.862.28'
    'E16' = '  +1.97%  '
    'D17' = '62.060.27'
    'E17' = '  +1.32%  '
    'D18' = '2.419.59'
    'E18' = '  +1.76%  '
    'D19' = '11.33'
    'E19' = '  +3.43%  '
    'E20' = '  +1.05%  '
    'D21' = '324.48'
    'E21' = '  +1.19%  '
    'D22' = '6.75'
    'E22' = '  +0.92%  '
    'E23' = '  +0.02%  '
    'E24' = '  +1.90%  '
    'E25' = '  -1.57%  '
    'D26' = '8.93'
    'E26' = '  +0.66%  '
    'D27' = '582.69'
    'E27' = '  +11.68%  '
    'D28' = '2.539.28'
    'E28' = '  +1.74%  '
    'D29' = '1.00'
    'E29' = '  -0.07%  '
    'D30' = '0.0₃0943'
    'E30' = '  +4.29%  '
    'E31' = '  +4.80%  '
    'D32' = '8.27'
    'E32' = '  +1.06%  '
    'E33' = '  +1.72%  '
    'E34' = '  +2.62%  '
    'E35' = '  +2.10%  '
    'D36' = '5.73'
    'E36' = '  +3.41%  '
    'E37' = '  +0.08%  '
    'D38' = '4.79'
    'E38' = '  +1.99%  '
    'D39' = '0.384'
    'E39' = '  +1.57%  '
    'D40' = '153.21'
    'E40' = '  +4.55%  '
    'D41' = '18.66'
    'E41' = '  +0.75%  '
    'E42' = '  -2.88%  '
    'E43' = '  -0.13%  '
    'E44' = '  +8.35%  '
    'D45' = '149.96'
    'E45' = '  +1.75%  '
    'E46' = '  +1.63%  '
    'E47' = '  +3.03%  '
    'D48' = '20.28'
    'E48' = '  +2.39%  '
    'E49' = '  +1.94%  '
    'E50' = '  +1.59%  '
    'E51' = '  +1.96%  '
}

# The Price column (D) holds values that look numeric (e.g. "561.68", "1.00"),
# but must stay plain text (matches the sheet's existing inline-string cells,
# preserves exact formatting like trailing zeros / thousands separators).
# Force the whole column to Text format before writing so Excel doesn't
# silently coerce the strings to numbers, then restore the default style
# afterwards so no visible formatting changes are introduced.
$ws.Range("D2:D51").NumberFormat = "@"

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$ws.Range("D2:D51").Style = "Normal"
